$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header labels from "darkNumber" to "reporting rate"
$ws.Range("C1").Value = "reporting rate (value)"
$ws.Range("D1").Value = "reporting rate (start index)"

# Remove the test data row (Luxembourg), shifting rows 3-5 up to 2-4
$ws.Rows(2).Delete()

# Update the selection to match the final state
$ws.Range("D9").Select()
